$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header
$ws.Range("D1").Value = "Full %"
$ws.Cells.Item(1, 4).Font.Bold = $true

# Data values for column D (rows 2-21)
$values = @(45, 62, 18, 73, 100, 63, 14, 27, 36, 94, 87, 30, 26, 43, 61, 47, 73, 25, 10, 0)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $values[$i]
}

$ws.Range("D21").Select()
